$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price snapshot refresh (GitHub Actions scheduled update).
# Column D holds prices as literal text (values like "70.536.29" use
# dotted thousands-grouping, not a single decimal point), so for the cells
# whose price looks like an ordinary decimal number (e.g. "4.84") we force
# the Text number format first -- otherwise Excel would helpfully (and
# wrongly) store it as a numeric value instead of the literal string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.614.08"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.556.32"
$ws.Range("E3").Value = "  -4.77%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.58"
$ws.Range("E5").Value = "  -3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.17"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.557.75"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.84"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.024.52"
$ws.Range("E14").Value = "  -4.78%  "
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.409.83"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.14"
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.611.67"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.74"
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.67"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  -8.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.99"
$ws.Range("E22").Value = "  -5.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.11"
$ws.Range("E26").Value = "  -5.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("E27").Value = "  -3.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.717.41"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.95"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "487.26"
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.04"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("E37").Value = "  +6.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.74"
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.85"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.78"
$ws.Range("E43").Value = "  -4.37%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "145.44"
$ws.Range("E47").Value = "  -7.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.57"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.533"
$ws.Range("E49").Value = "  -5.93%  "
$ws.Range("E50").Value = "  -6.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.596"
$ws.Range("E51").Value = "  -1.67%  "
